$d = $word.ActiveDocument

# Renaming an inline picture's docPr/name in a footer directly via
# "<HeaderFooter>.Range.InlineShapes(n)" can trip a stale-handle error in
# this host, so resolve the shape through the paragraph that actually
# contains it (Range.Paragraphs(i).Range.InlineShapes(1)) instead - that
# addressing path is stable for both headers and footers.
function Get-ShapeInRange($range) {
    $count = $range.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $paraRange = $range.Paragraphs($i).Range
        if ($paraRange.InlineShapes.Count -gt 0) {
            return $paraRange.InlineShapes(1)
        }
    }
    return $null
}

# footer1.xml (w:type="first" footer, docPr id="3") - PearsonLogo.png:
# image2.png -> image1.png
$footerFirst = $d.Sections(1).Footers(2)
$shape = Get-ShapeInRange $footerFirst.Range
$shape.Name = "image1.png"

# footer2.xml (w:type="default" footer, docPr id="2") - PearsonLogo.png:
# image2.png -> image1.png
$footerDefault = $d.Sections(1).Footers(1)
$shape = Get-ShapeInRange $footerDefault.Range
$shape.Name = "image1.png"

# header1.xml (w:type="first" header, docPr id="1") - BTec_Logo-Orange:
# image1.jpg -> image2.jpg
$headerFirst = $d.Sections(1).Headers(2)
$shape = Get-ShapeInRange $headerFirst.Range
$shape.Name = "image2.jpg"
